$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "v"
$ws.Range("E6").Value = "v"

$ws.Range("E8").Select()
